# Fruta / hortaliza, semanal
# Insert a new weekly record before the existing row 23, which pushes
# rows 23:48 down to 24:49 (and extends the used range to A1:T49).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 23 - this shifts rows 23-48 down to 24-49,
# carrying all their existing data and formatting with them.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly record.
$ws.Cells.Item(23, 1).Value = 6
$ws.Cells.Item(23, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(23, 3).Value = "Metropolitana"
$ws.Cells.Item(23, 4).Value = 44544
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100101
$ws.Cells.Item(23, 8).Value = "Berries"
$ws.Cells.Item(23, 9).Value = 100101008
$ws.Cells.Item(23, 10).Value = "Mora"
$ws.Cells.Item(23, 11).Value = "Sin especificar"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 250
$ws.Cells.Item(23, 14).Value = 5000
$ws.Cells.Item(23, 15).Value = 5000
$ws.Cells.Item(23, 16).Value = 5000
$ws.Cells.Item(23, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(23, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(23, 19).Value = 2500
$ws.Cells.Item(23, 20).Value = 2
